$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell values that changed
$ws.Range("L4").Value = 659773
$ws.Range("B6").Value = 22482006
$ws.Range("L7").Value = 134036
$ws.Range("B10").Value = 22482007

# Update the selection on the sheet (was L7, now L2:L10 with active cell L2)
$ws.Range("L2:L10").Select()

# Update the workbook window position/size
$excel.Windows.Item(1).Left = 38420
$excel.Windows.Item(1).Top = 1640
$excel.Windows.Item(1).Width = 21600
$excel.Windows.Item(1).Height = 20180
